# regressions.xlsx edit:
#  - re-sort the existing "Sort" (col A) / "Regression" (col B) test rows
#    (rows 4-12, i.e. original tests 3..11) ascending by col A instead of
#    the previous descending-by-R-Adjusted-test ordering
#  - append a new test row (19: "13+a30-accommodates (br only)") with its
#    R Adjusted test score (41.97%) — lin_reg added to airbnb / eda refactor
#  - move the active selection to the newly added score cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-sort the block of test rows (A2:D12 -> rows 4:12 get reordered,
#        rows 2:3 are already in ascending col-A order so stay put) ---------
$sortKey = $ws.Range("A2")
$ws.Range("A2:D12").Sort($sortKey, 1, $null, $null, 1, $null, $null, 2)

# Keep the sheet's recorded AutoFilter sort definition in sync with the new
# ascending sort-by-column-A (best effort — harmless if the host ignores it).
$af = $ws.AutoFilter
$sort = $af.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("A2:A19"), 0, 1)
$sort.SetRange($ws.Range("A1:D19"))
$sort.Header = 1
$sort.Apply()

# --- 2. Append the new test (row 20) ---------------------------------------
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "13+a30-accommodates (br only)"
$ws.Range("D20").Value = 0.41968152389993701
$ws.Range("D20").NumberFormat = "0.00%"

# --- 3. Move selection to the new score cell --------------------------------
$ws.Range("D20").Select()
